# Apply updated cryptocurrency price/volume figures to the sheet.
# Generated from the commit diff: updates columns D (Price) and E (Volume(1h))
# for the rows that changed. D-column values that look like plain numbers are
# prefixed with a leading apostrophe so Excel keeps them as text (matching the
# original inlineStr/text storage) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.652.67"
$ws.Range("E2").Value = "  +5.67%  "
$ws.Range("D3").Value = "3.187.45"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'401.98"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("D6").Value = "'108.58"
$ws.Range("E6").Value = "  +5.01%  "
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +5.42%  "
$ws.Range("D10").Value = "'39.04"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "3.682.15"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").Value = "'19.04"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("E15").Value = "  +3.44%  "
$ws.Range("E16").Value = "  +8.83%  "
$ws.Range("D17").Value = "3.187.61"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "'10.57"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").Value = "54.545.60"
$ws.Range("E19").Value = "  +5.20%  "
$ws.Range("E20").Value = "  +3.80%  "
$ws.Range("E21").Value = "  +3.59%  "
$ws.Range("D22").Value = "0.0₃0999"
$ws.Range("E22").Value = "  +3.13%  "
$ws.Range("D23").Value = "'72.76"
$ws.Range("E23").Value = "  +3.89%  "
$ws.Range("D24").Value = "'275.70"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("E25").Value = "  +5.05%  "
$ws.Range("D26").Value = "'8.08"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "'27.79"
$ws.Range("E27").Value = "  +2.52%  "
$ws.Range("D28").Value = "'7.47"
$ws.Range("E28").Value = "  +3.71%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  +3.94%  "
$ws.Range("E32").Value = "  +6.76%  "
$ws.Range("D33").Value = "'0.0504"
$ws.Range("E33").Value = "  +12.37%  "
$ws.Range("D34").Value = "'36.94"
$ws.Range("E34").Value = "  +3.45%  "
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").Value = "'3.65"
$ws.Range("E37").Value = "  +7.38%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  +10.05%  "
$ws.Range("E40").Value = "  +13.09%  "
$ws.Range("E41").Value = "  +3.39%  "
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").Value = "'17.31"
$ws.Range("E43").Value = "  +1.93%  "
$ws.Range("D44").Value = "'130.86"
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").Value = "'22.34"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").Value = "'2.06"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").Value = "2.091.20"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("E50").Value = "  +9.71%  "
$ws.Range("E51").Value = "  +12.27%  "
